# ProjectTimeline.docx edit:
#  1. Append a new table row documenting the equator-system rework.
#  2. Split the run "with the perlin temperature system." (in the row
#     that was previously last) so that the lower-case "p" of "perlin"
#     becomes an upper-case "P" living in its own run, matching the
#     reworked OOXML from the diff.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Remember which row currently holds the "...perlin temperature system."
# text before we append anything after it.
$notesRowIndex = $t.Rows.Count

# --- Part 1: add the new row describing the equator rework ---
$newRow = $t.Rows.Add()

$newRow.Cells.Item(1).Range.Text = "28/11/2021"
$newRow.Cells.Item(2).Range.Text = "1 Hour 50 minutes"
$newRow.Cells.Item(3).Range.Text = "World Generation – Objective 1"
$newRow.Cells.Item(4).Range.Text = "Reworked the equator system to blend better with the surroundings, also removed a significant amount of the randomness of the system to provide a more consistently positive result."

# --- Part 2: capitalise "perlin" -> "Perlin", leaving it split across runs ---
$notesCell = $t.Rows.Item($notesRowIndex).Cells.Item(4)
$findRange = $notesCell.Range.Duplicate
$found = $findRange.Find.Execute("perlin", $true, $true, $false, $false, $false, $true)

if ($found) {
    $letterStart = $findRange.Start
    $pRange = $d.Range($letterStart, $letterStart + 1)
    # Toggling a character property and clearing it again keeps this run
    # distinct from its neighbours instead of letting it re-merge with
    # them once the identical text formatting is restored.
    $pRange.Bold = 1
    $pRange.Text = "P"
    $pRange2 = $d.Range($letterStart, $letterStart + 1)
    $pRange2.Bold = 0
}
